$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date for every existing data row (2-27)
# from 2024-02-16 (45338) to 2024-04-08 (45390).
$ws.Range("C2:C27").Value = 45390

# Row 27 picks up an explicit row height in the new file.
$ws.Rows.Item(27).RowHeight = 15

# Append the new record as row 28.
$ws.Range("A28").Value = "A 11790-2024"

$ws.Range("B28").Value = 45374
$ws.Range("B28").NumberFormat = "YYYY-MM-DD"

$ws.Range("C28").Value = 45390
$ws.Range("C28").NumberFormat = "YYYY-MM-DD"

$ws.Range("D28").Value = "OKÄNT"
$ws.Range("E28").Value = "OKÄNT"

$ws.Range("G28").Value = 0.7
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0
$ws.Range("N28").Value = 0
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = 0
$ws.Range("Q28").Value = 0

# R28 stays blank but carries the wrap-text style used throughout column R.
$ws.Range("R28").WrapText = $true
